# Generate Report for handoff
# The f12bf0be-48f2-4685-b862-21f679bd3c82 file has been handed off again,
# so its status moves from "Handed back: in sync with en-US" to
# "Ready for handoff", and the "Latest Handoff Datetime" for both files
# is refreshed to reflect the new handoff run.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
# Row 3 corresponds to f12bf0be-48f2-4685-b862-21f679bd3c82.md
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet ------------------------------------------------------
# Row 2 -> 29a97ec9-0cb8-4386-afa4-0fc85831dfb1 ; refresh handoff datetime
$zhcn.Range("D2").Value = "2016-02-16 10:38:16"

# Row 3 -> f12bf0be-48f2-4685-b862-21f679bd3c82 ; now ready for handoff again
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-02-16 10:38:16"

# --- de-de sheet ------------------------------------------------------
# Row 2 -> 29a97ec9-0cb8-4386-afa4-0fc85831dfb1 ; refresh handoff datetime
$dede.Range("D2").Value = "2016-02-16 10:38:30"

# Row 3 -> f12bf0be-48f2-4685-b862-21f679bd3c82 ; now ready for handoff again
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-02-16 10:38:30"
